# Insert a new weekly record for "Apio" (Macroferia Regional de Talca) as row 295,
# shifting the existing rows 295-319 down to 296-320 (dimension grows to A1:R320).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 295..319 down by one to make room for the new record.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new observation.
$ws.Range("A295").Value = 5
$ws.Range("B295").Value = "Macroferia Regional de Talca"
$ws.Range("C295").Value = "Maule"
$ws.Range("D295").Value = 45106
$ws.Range("E295").Value = 7
$ws.Range("F295").Value = 100112017
$ws.Range("G295").Value = "Apio"
$ws.Range("H295").Value = "Americana (o)"
$ws.Range("I295").Value = "Primera"
$ws.Range("J295").Value = 700
$ws.Range("K295").Value = 6000
$ws.Range("L295").Value = 6000
$ws.Range("M295").Value = 6000
$ws.Range("N295").Value = "`$/docena de matas"
$ws.Range("O295").Value = "Provincia del Elquí"
$ws.Range("P295").Value = 1000
$ws.Range("Q295").Value = 6
$ws.Range("R295").Value = "Hortaliza"
